$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <-> Row 3 swap (title + uri)
$ws.Range("A2").Value = "Lee County, Alabama"
$ws.Range("E2").Value = "http://www.tornadohistoryproject.com/tornado/Alabama/Lee/map"
$ws.Range("A3").Value = "NOAA National Weather Service"
$ws.Range("E3").Value = "https://www.weather.gov/chs/LibertyCountytornado2019"

# Rows 4,5,6 rotate: new4=old5, new5=old6, new6=old4
$ws.Range("A4").Value = "County Road 79 Tornado - March 3, 2019"
$ws.Range("E4").Value = "https://www.weather.gov/bmx/event_03032019cr79"
$ws.Range("A5").Value = "Tornadoes of March 3, 2019"
$ws.Range("E5").Value = "https://www.weather.gov/bmx/event_03032019"
$ws.Range("A6").Value = "Davisville-Corbett Crossroad Tornado - March 3, 2019"
$ws.Range("E6").Value = "https://www.weather.gov/bmx/event_03032019davisville"
